$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename existing sheet and add the new "PassengerDetails" sheet right
#    after it (so order becomes: LoginDetails, PassengerDetails).
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "PassengerDetails"
$ws1.Name = "LoginDetails"

# ---------------------------------------------------------------------------
# 2. PassengerDetails sheet - header row
# ---------------------------------------------------------------------------
$ws2.Range("A1").Value = "Username"
$ws2.Range("B1").Value = "Password"
$ws2.Range("C1").Value = "FirstName"
$ws2.Range("D1").Value = "LastName"
$ws2.Range("E1").Value = "Date"
$ws2.Range("F1").Value = "Month"
$ws2.Range("G1").Value = "Year"
$ws2.Range("H1").Value = "Country"
$ws2.Range("I1").Value = "Address"
$ws2.Range("J1").Value = "City"
$ws2.Range("K1").Value = "Mnumber"
$ws2.Range("L1").Value = "Email"
$ws2.Range("M1").Value = "CardNumber"
$ws2.Range("N1").Value = "CardName"
$ws2.Range("O1").Value = "CardMonth"
$ws2.Range("P1").Value = "CardYear"
$ws2.Range("Q1").Value = "CVV"

# ---------------------------------------------------------------------------
# 3. PassengerDetails sheet - data row 2
# ---------------------------------------------------------------------------
$ws2.Range("C2").Value = "Tushar"
$ws2.Range("D2").Value = "Jadhav"
$ws2.Range("E2").Value = 3
$ws2.Range("F2").Value = 8
$ws2.Range("G2").Value = 1999
$ws2.Range("H2").Value = "India"
$ws2.Range("I2").Value = "Asalfa"
$ws2.Range("J2").Value = "Mumbai"
$ws2.Range("K2").Value = 8745963287
$ws2.Range("L2").Value = "tusharjadhav123@gmail.com"
$ws2.Range("M2").Value = "4111 1111 1111 1111"
$ws2.Range("N2").Value = "Tushar Jadhav"
$ws2.Range("O2").NumberFormat = "@"
$ws2.Range("O2").Value = "08"
$ws2.Range("P2").Value = 2026
$ws2.Range("Q2").Value = 123

# ---------------------------------------------------------------------------
# 4. PassengerDetails sheet - data row 3
# ---------------------------------------------------------------------------
$ws2.Range("C3").Value = "Anjali"
$ws2.Range("D3").Value = "Patil"
$ws2.Range("E3").Value = 17
$ws2.Range("F3").Value = 11
$ws2.Range("G3").Value = 2001
$ws2.Range("H3").Value = "India"
$ws2.Range("I3").Value = "Asalfa"
$ws2.Range("J3").Value = "Mumbai"
$ws2.Range("K3").Value = 8745963287
$ws2.Range("L3").Value = "tusharjadhav123@gmail.com"
$ws2.Range("M3").Value = "4111 1111 1111 1111"
$ws2.Range("N3").Value = "Tushar Jadhav"
$ws2.Range("O3").NumberFormat = "@"
$ws2.Range("O3").Value = "08"
$ws2.Range("P3").Value = 2028
$ws2.Range("Q3").Value = 456

# ---------------------------------------------------------------------------
# 5. PassengerDetails sheet - A2/B2/A3/B3 carry the login e-mail/password
#    values together with mailto hyperlinks (mirroring LoginDetails).
# ---------------------------------------------------------------------------
$ws2.Range("A2").Value = "tusharjadhav228@gmail.com"
$ws2.Hyperlinks.Add($ws2.Range("A2"), "mailto:tusharjadhav228@gmail.com") | Out-Null
$ws2.Range("A2").Style = "Hyperlink"

$ws2.Range("B2").Value = "Sakshi@1228"
$ws2.Hyperlinks.Add($ws2.Range("B2"), "mailto:Sakshi@1228") | Out-Null
$ws2.Range("B2").Style = "Hyperlink"

$ws2.Range("A3").Value = "tushar.jadhav.work@gmail.com"
$ws2.Hyperlinks.Add($ws2.Range("A3"), "mailto:tushar.jadhav.work@gmail.com") | Out-Null
$ws2.Range("A3").Style = "Hyperlink"

$ws2.Range("B3").Value = "Jadhav@1228"
$ws2.Hyperlinks.Add($ws2.Range("B3"), "mailto:Jadhav@1228") | Out-Null
$ws2.Range("B3").Style = "Hyperlink"

# ---------------------------------------------------------------------------
# 6. PassengerDetails sheet - column widths (best match achievable)
# ---------------------------------------------------------------------------
$ws2.Columns.Item(1).ColumnWidth = 27
$ws2.Columns.Item(2).ColumnWidth = 16.333333333333332
$ws2.Columns.Item(3).ColumnWidth = 9
$ws2.Columns.Item(4).ColumnWidth = 9.333333333333334
$ws2.Columns.Item(5).ColumnWidth = 4.833333333333333
$ws2.Columns.Item(8).ColumnWidth = 9.166666666666666
$ws2.Columns.Item(9).ColumnWidth = 10.333333333333334
$ws2.Columns.Item(11).ColumnWidth = 10.666666666666666
$ws2.Columns.Item(12).ColumnWidth = 25
$ws2.Columns.Item(13).ColumnWidth = 17
$ws2.Columns.Item(14).ColumnWidth = 12
$ws2.Columns.Item(15).ColumnWidth = 10.333333333333334

# ---------------------------------------------------------------------------
# 7. Selections / active tab - LoginDetails keeps A3:B3 selected (no longer
#    the visible tab), PassengerDetails becomes the active tab with N9
#    selected.
# ---------------------------------------------------------------------------
$ws1.Range("A3:B3").Select() | Out-Null
$ws2.Activate()
$ws2.Range("N9").Select() | Out-Null
